$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 9 (shifts old rows 9-17 down to 10-18)
$ws.Rows("9:9").Insert()

# Copy formatting from row 10 (the row that used to be row 9) onto the new row 9
$ws.Range("B10:L10").Copy()
$ws.Range("B9:L9").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row 9 content ("exposure_time" / Timer1 entry)
$ws.Range("B9").Value2 = "uint8      'E'"
$ws.Range("C9").Value2 = "uint16                       exposure_time"
$ws.Range("L9").Value2 = "struct Timer1 T1"

# Merge C9:D9 like the row below it
$ws.Range("C9:D9").Merge()

# Select the newly merged cell, matching the final user action
$ws.Range("C9:D9").Select()
